$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - "Save", styled like the other header cells (copy style from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for the new "Save" column (H2:H14)
$saveValues = @(1, 0, 1, 0, 1, 1, 1, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
